$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range('D2')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '64.743.30'
$rng.Style = $origStyle
$ws.Range('E2').Value = '  +1.46%  '
$rng = $ws.Range('D3')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '3.453.52'
$rng.Style = $origStyle
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.00%  '
$rng = $ws.Range('D5')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '575.01'
$rng.Style = $origStyle
$ws.Range('E5').Value = '  +0.69%  '
$rng = $ws.Range('D6')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '160.85'
$rng.Style = $origStyle
$ws.Range('E6').Value = '  +2.25%  '
$rng = $ws.Range('D7')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.615'
$rng.Style = $origStyle
$ws.Range('E7').Value = '  +12.58%  '
$ws.Range('E8').Value = '  +0.06%  '
$rng = $ws.Range('D9')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '3.454.37'
$rng.Style = $origStyle
$ws.Range('E9').Value = '  +0.95%  '
$rng = $ws.Range('D10')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '7.21'
$rng.Style = $origStyle
$ws.Range('E10').Value = '  -2.51%  '
$ws.Range('E11').Value = '  +1.59%  '
$rng = $ws.Range('D12')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.446'
$rng.Style = $origStyle
$ws.Range('E12').Value = '  +3.17%  '
$rng = $ws.Range('D13')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '4.042.28'
$rng.Style = $origStyle
$ws.Range('E13').Value = '  +0.74%  '
$rng = $ws.Range('D14')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.135'
$rng.Style = $origStyle
$ws.Range('E14').Value = '  +0.35%  '
$rng = $ws.Range('D15')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0000193'
$rng.Style = $origStyle
$ws.Range('E15').Value = '  -0.45%  '
$rng = $ws.Range('D16')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '28.17'
$rng.Style = $origStyle
$ws.Range('E16').Value = '  +3.40%  '
$rng = $ws.Range('D17')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '64.841.57'
$rng.Style = $origStyle
$ws.Range('E17').Value = '  +1.71%  '
$rng = $ws.Range('D18')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '3.466.08'
$rng.Style = $origStyle
$ws.Range('E18').Value = '  -0.04%  '
$rng = $ws.Range('D19')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '6.44'
$rng.Style = $origStyle
$ws.Range('E19').Value = '  +2.39%  '
$rng = $ws.Range('D20')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '14.32'
$rng.Style = $origStyle
$ws.Range('E20').Value = '  +1.55%  '
$rng = $ws.Range('D21')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '379.69'
$rng.Style = $origStyle
$ws.Range('E21').Value = '  -0.41%  '
$rng = $ws.Range('D22')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '8.10'
$rng.Style = $origStyle
$ws.Range('E22').Value = '  +0.12%  '
$rng = $ws.Range('D23')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.549'
$rng.Style = $origStyle
$ws.Range('E23').Value = '  +3.48%  '
$rng = $ws.Range('D24')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '1.00'
$rng.Style = $origStyle
$ws.Range('E24').Value = '  +0.08%  '
$rng = $ws.Range('D25')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '72.18'
$rng.Style = $origStyle
$ws.Range('E25').Value = '  +0.49%  '
$rng = $ws.Range('D26')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0000120'
$rng.Style = $origStyle
$ws.Range('E26').Value = '  -1.53%  '
$rng = $ws.Range('D27')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '10.04'
$rng.Style = $origStyle
$ws.Range('E27').Value = '  +7.00%  '
$ws.Range('E28').Value = '  -0.39%  '
$rng = $ws.Range('D29')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '1.00'
$rng.Style = $origStyle
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('E30').Value = '  +9.10%  '
$rng = $ws.Range('D31')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '6.13'
$rng.Style = $origStyle
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('E32').Value = '  +0.95%  '
$rng = $ws.Range('D33')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '23.58'
$rng.Style = $origStyle
$ws.Range('E33').Value = '  +1.20%  '
$rng = $ws.Range('D34')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '7.12'
$rng.Style = $origStyle
$ws.Range('E34').Value = '  +4.73%  '
$rng = $ws.Range('D35')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '1.62'
$rng.Style = $origStyle
$ws.Range('E35').Value = '  +11.84%  '
$rng = $ws.Range('D36')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '161.48'
$rng.Style = $origStyle
$ws.Range('E36').Value = '  +1.08%  '
$rng = $ws.Range('D37')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '1.92'
$rng.Style = $origStyle
$ws.Range('E37').Value = '  +5.17%  '
$rng = $ws.Range('D38')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0775'
$rng.Style = $origStyle
$ws.Range('E38').Value = '  +2.39%  '
$rng = $ws.Range('D39')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '2.961.88'
$rng.Style = $origStyle
$ws.Range('E39').Value = '  -0.79%  '
$rng = $ws.Range('D40')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '26.48'
$rng.Style = $origStyle
$ws.Range('E40').Value = '  -2.25%  '
$ws.Range('E45').Value = '  +1.60%  '
$rng = $ws.Range('D46')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '25.60'
$rng.Style = $origStyle
$ws.Range('E46').Value = '  +10.13%  '
$ws.Range('E47').Value = '  +1.99%  '
$ws.Range('E48').Value = '  +8.35%  '
$rng = $ws.Range('D49')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '311.59'
$rng.Style = $origStyle
$ws.Range('E49').Value = '  +6.42%  '
$rng = $ws.Range('D50')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '6.62'
$rng.Style = $origStyle
$ws.Range('E50').Value = '  +4.18%  '
$rng = $ws.Range('D51')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.863'
$rng.Style = $origStyle
$ws.Range('E51').Value = '  +3.43%  '

# Row swaps: Filecoin/RenderToken (41<->42) and VeChain/OKB (43<->44)

$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$rng = $ws.Range('D41')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '4.56'
$rng.Style = $origStyle
$ws.Range('E41').Value = '  +5.59%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$rng = $ws.Range('D42')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '6.58'
$rng.Style = $origStyle
$ws.Range('E42').Value = '  +3.44%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$rng = $ws.Range('D43')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0317'
$rng.Style = $origStyle
$ws.Range('E43').Value = '  +0.52%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$rng = $ws.Range('D44')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '42.75'
$rng.Style = $origStyle
$ws.Range('E44').Value = '  +1.75%  '

Write-Host "Done"
